# Updated cryptos list with GitHub Actions
# Applies the diff: refresh Price (column D) and Volume(1h) (column E)
# values for the crypto rows on Sheet1.
#
# Cells whose new text is a plain numeric-looking string (e.g. "1.00",
# "29.00", "41.20") are explicitly formatted as Text ("@") before the
# assignment so Excel's COM layer keeps them as literal strings (matching
# the source workbook's original inline-string cells) instead of silently
# coercing them to numbers and dropping formatting such as trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.367.01"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "3.080.88"
$ws.Range("E3").Value = "  +3.90%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.29"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.13"
$ws.Range("E6").Value = "  +3.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.209"
$ws.Range("E9").Value = "  +5.95%  "
$ws.Range("D10").Value = "3.083.05"
$ws.Range("E10").Value = "  +4.03%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.21"
$ws.Range("E13").Value = "  +6.07%  "
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.00"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "76.375.44"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("D18").Value = "3.069.83"
$ws.Range("E18").Value = "  +3.70%  "
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.16"
$ws.Range("E20").Value = "  +5.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.30"
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("E22").Value = "  +9.27%  "
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("D24").Value = "3.242.85"
$ws.Range("E24").Value = "  +4.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.34"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.35"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.29"
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "501.79"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.91"
$ws.Range("E34").Value = "  +4.35%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.124"
$ws.Range("E36").Value = "  +12.51%  "
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.70"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.05"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "193.22"
$ws.Range("E40").Value = "  +7.04%  "
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("E42").Value = "  -9.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.795"
$ws.Range("E44").Value = "  +20.37%  "
$ws.Range("E45").Value = "  +3.69%  "
$ws.Range("E46").Value = "  +4.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.20"
$ws.Range("E47").Value = "  +2.58%  "
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  +5.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.597"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("E51").Value = "  -0.55%  "
